$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '26.077.93'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -0.52%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.660.50'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -1.21%  '
$ws.Range("E4").Value = '  -0.18%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '207.46'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -1.84%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.5154'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -2.83%  '
$ws.Range("E7").Value = '  -0.17%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.2578'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -3.92%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06278'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -0.47%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '20.93'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -2.02%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07514'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -0.18%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.635.83'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -2.88%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '4.398'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -1.85%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.5378'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -5.17%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '66.14'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -0.61%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.0₅7901'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -2.93%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '26.103.27'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -0.64%  '
$ws.Range("E18").Value = '  -0.16%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '4.690'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -3.44%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '186.86'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -1.03%  '
$ws.Range("E21").Value = '  -3.68%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.171'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -0.77%  '
$ws.Range("E23").Value = '  -0.09%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '148.30'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +0.49%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.1209'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -4.02%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '7.377'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -3.29%  '
$ws.Range("E27").Value = '  -1.93%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.375'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +2.15%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.06114'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -5.19%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.261'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -2.04%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '3.465'
$ws.Range("D31").Style = "Normal"
$ws.Range("E32").Value = '  -2.64%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.627'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -1.77%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.9836'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -2.68%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.388'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -1.19%  '
$ws.Range("E36").Value = '  +1.06%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.5866'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -4.02%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.105.83'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +0.37%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.01589'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -1.80%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '5.967'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -3.41%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.8460'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -2.54%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.003'
$ws.Range("D42").Style = "Normal"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '100.00'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -0.32%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.809.78'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -1.21%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0₈107'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -3.62%  '
$ws.Range("E46").Value = '  +0.02%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '54.83'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -3.85%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '7.991'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -0.04%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.05231'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -0.67%  '
$ws.Range("E50").Value = '  -0.68%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '5.849'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -2.02%  '
